# Re-style the three data tables (slides 14, 15, 16) from the deck's
# custom "Table_0" style to PowerPoint's built-in "No Style, Table Grid"
# style, matching what the Table Design gallery would apply.

$p = $ppt.ActivePresentation

$newStyleId = "{7A190B7B-183B-4223-912F-E5C288AA35B5}"
$slideIndexesWithTables = @(14, 15, 16)

foreach ($slideIndex in $slideIndexesWithTables) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
